# Kapitel 5 & 6 getauscht, Zuordnungsdokument angepasst
#
# "5. Bildakquise und Datenaufbereitung " (Tobias Rempel, rows 29-45) and
# "6. SonoScape Analyse " (Andy Kruder, rows 46-49) swap places and are
# renumbered: SonoScape Analyse becomes chapter 5 (rows 29-32) and
# Bildakquise und Datenaufbereitung becomes chapter 6 (rows 33-49).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New chapter 5: SonoScape Analyse (formerly chapter 6, Andy Kruder)
$ws.Range("A29").Value = "5. SonoScape Analyse "
$ws.Range("B29").Value = "Andy Kruder"
$ws.Range("C29").Value = $null

$ws.Range("A30").Value = "5.1. System Struktur "
$ws.Range("B30").Value = "Andy Kruder"
$ws.Range("C30").Value = $null

$ws.Range("A31").Value = "5.2. Root Zugriff "
$ws.Range("B31").Value = "Andy Kruder"
$ws.Range("C31").Value = $null

$ws.Range("A32").Value = "5.3. Netzwerk Verbindung"
$ws.Range("B32").Value = "Andy Kruder"
$ws.Range("C32").Value = $null

# New chapter 6: Bildakquise und Datenaufbereitung (formerly chapter 5, Tobias Rempel)
$ws.Range("A33").Value = "6. Bildakquise und Datenaufbereitung "
$ws.Range("B33").Value = "Tobias Rempel"
$ws.Range("C33").Value = 7080879

$ws.Range("A34").Value = "6.1. Bildakquise "
$ws.Range("B34").Value = "Tobias Rempel"
$ws.Range("C34").Value = 7080879

$ws.Range("A35").Value = "6.2. Kodierung "
$ws.Range("B35").Value = "Tobias Rempel"
$ws.Range("C35").Value = 7080879

$ws.Range("A36").Value = "6.3. Dekodierung"
$ws.Range("B36").Value = "Tobias Rempel"
$ws.Range("C36").Value = 7080879

$ws.Range("A37").Value = "6.4. Implementierung mit FFmpeg "
$ws.Range("B37").Value = "Tobias Rempel"
$ws.Range("C37").Value = 7080879

$ws.Range("A38").Value = "6.4.1. Auswahl eines Videocodecs "
$ws.Range("B38").Value = "Tobias Rempel"
$ws.Range("C38").Value = 7080879

$ws.Range("A39").Value = "6.4.2. Nutzung des MediaCodec "
$ws.Range("B39").Value = "Tobias Rempel"
$ws.Range("C39").Value = 7080879

$ws.Range("A40").Value = "6.4.3. Verarbeitung des Datenstroms"
$ws.Range("B40").Value = "Tobias Rempel"
$ws.Range("C40").Value = 7080879

$ws.Range("A41").Value = "6.4.4. Konfiguration des MediaCodec "
$ws.Range("B41").Value = "Tobias Rempel"
$ws.Range("C41").Value = 7080879

$ws.Range("A42").Value = "6.4.5. Performance und Funktionalität "
$ws.Range("B42").Value = "Tobias Rempel"
$ws.Range("C42").Value = 7080879

$ws.Range("A43").Value = "6.5. Implementierung mit X11Lib und LZ4 "
$ws.Range("B43").Value = "Tobias Rempel"
$ws.Range("C43").Value = 7080879

$ws.Range("A44").Value = "6.5.1. Bildakquise mit X11Lib "
$ws.Range("B44").Value = "Tobias Rempel"
$ws.Range("C44").Value = 7080879

$ws.Range("A45").Value = "6.5.2. Datenreduktion "
$ws.Range("B45").Value = "Tobias Rempel"
$ws.Range("C45").Value = 7080879

$ws.Range("A46").Value = "6.5.3. Kompression "
$ws.Range("B46").Value = "Tobias Rempel"
$ws.Range("C46").Value = 7080879

$ws.Range("A47").Value = "6.5.4. Verarbeitung auf dem Smartphone "
$ws.Range("B47").Value = "Tobias Rempel"
$ws.Range("C47").Value = 7080879

$ws.Range("A48").Value = "6.5.5. Performance und Funktionalität "
$ws.Range("B48").Value = "Tobias Rempel"
$ws.Range("C48").Value = 7080879

$ws.Range("A49").Value = "6.5.6. Optimierungspotenzial "
$ws.Range("B49").Value = "Tobias Rempel"
$ws.Range("C49").Value = 7080879

# Update the view state to match the author's final scroll/selection position
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 22
$aw.ScrollColumn = 1
$ws.Range("A56").Select()
